# Apply the APP-000017 assets_liabilities.xlsx update:
#  - Summary sheet: new client name + refreshed income/asset/liability/net-worth/ratio figures
#  - Assets sheet: drop the "Vehicles / Premium Car" line item, keep Liquid Assets + recompute TOTAL ASSETS
#  - Liabilities sheet: refreshed Credit Card balance + monthly payment and TOTAL LIABILITIES

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Cells.Item(3, 2).Value = "Anwar Al Shehhi"   # B3 Name
$summary.Cells.Item(4, 2).Value = 8298.18             # B4 Monthly Income (AED)
$summary.Cells.Item(6, 2).Value = 7765                # B6 Total Assets (AED)
$summary.Cells.Item(7, 2).Value = 10319               # B7 Total Liabilities (AED)
$summary.Cells.Item(8, 2).Value = -2554               # B8 Net Worth (AED)
$summary.Cells.Item(9, 2).Value = 0.75                # B9 Asset/Liability Ratio

# ---------------------------------------------------------------------------
# Assets sheet - remove the "Vehicles" row entirely (row 2), which shifts
# "Liquid Assets" up to row 2 and "TOTAL ASSETS" up to row 3.
# ---------------------------------------------------------------------------
$assets = $wb.Worksheets.Item("Assets")
$assets.Rows.Item(2).Delete()

$assets.Cells.Item(2, 3).Value = 7765   # C2 Liquid Assets / Savings Account value
$assets.Cells.Item(3, 3).Value = 7765   # C3 TOTAL ASSETS

# ---------------------------------------------------------------------------
# Liabilities sheet
# ---------------------------------------------------------------------------
$liabilities = $wb.Worksheets.Item("Liabilities")
$liabilities.Cells.Item(2, 3).Value = 10319   # C2 Credit Card Balance amount
$liabilities.Cells.Item(2, 4).Value = 516     # D2 Monthly Payment
$liabilities.Cells.Item(3, 3).Value = 10319   # C3 TOTAL LIABILITIES
